$d = $word.ActiveDocument

# Remove the two empty paragraphs that originally separated the text blocks
# (two empty paragraphs still exist at the very end of the document afterwards).
$d.Paragraphs.Item(2).Range.Delete()   # empty paragraph that followed "Matricule..."
$d.Paragraphs.Item(3).Range.Delete()   # empty paragraph that followed "Quand..."

# --- Paragraph 1: "Matricule ..." -> new ending text ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$r1.Text = "Matricule doit etre string ou int : parce que si int 00012 devient 12 -> fait"

# --- Paragraph 2: "Quand ..." -> the trailing comment is replaced with a single space ---
$p2 = $d.Paragraphs.Item(2)
$find2 = $p2.Range.Find
$find2.Execute(" -> j’ai fait pour chantier", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null

# Paragraphs 3 "Mettre automatique..." and 4 "Ajouter marque..." stay unchanged.

# --- Insert the 8 new paragraphs after paragraph 4 ("Ajouter marque..."), before the trailing empty paragraph ---
$p4 = $d.Paragraphs.Item(4)
$insertPos = $p4.Range.End - 1
$ins = $d.Range($insertPos, $insertPos)
$newText = [char]13 + "Modifier les vue pour ne pas utilise dispose sur la vue" + [char]13 + `
           "Est ce qu’on laisse toutes les pop ups ? " + [char]13 + `
           "Est ce que la taille du nom d’un chantier dans livraison est assez grande ? " + [char]13 + `
           "Est ce qu’on supprime tous les objet qui vont avec un genre on supprime les livraison d’un chantier si on supprime le chantier ?? -> fait" + [char]13 + `
           "On ne peut pas utiliser un numero matricule de quelquun archive " + [char]13 + `
           "Numero de matricule supprime toujours les 0 -> fait" + [char]13 + `
           "Enlever el fait que le bouton employe est selectionne de base" + [char]13 + `
           "Mettre le truc archive et supprime ? "
$ins.InsertAfter($newText)

# --- The document ends with two empty paragraphs; add the extra one right before the
#     pre-existing trailing empty paragraph so that one of them stays untouched. ---
$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastStart = $lastP.Range.Start
$insEmpty = $d.Range($lastStart, $lastStart)
$insEmpty.InsertAfter([char]13)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
